# Append the professor's follow-up notes (a blank spacer line, a new date
# line, and a new "related work" bullet) to the end of the document, after
# the existing bulleted list.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1) A bare, completely empty spacer paragraph ----------------------
# (Using raw OOXML insertion here, rather than InsertParagraphAfter, keeps
# this paragraph genuinely empty -- no inherited list/style and no
# leftover empty run.)
$endRange = $d.Content
$endRange.Collapse(0)
$null = $endRange.InsertXML('<w:p ' + $wNs + '/>')

# --- 2) A new, plain date paragraph: "6/6/2022" -------------------------
$blankPara = $d.Paragraphs.Last
$blankPara.Range.InsertParagraphAfter()
$datePara = $d.Paragraphs.Last
$datePara.Range.Text = "6/6/2022"

# --- 3) A new top-level bullet in the same list (numId 1) as the other
#        bullets, continuing the existing "ListParagraph" bulleted list. --
$datePara.Range.InsertParagraphAfter()
$listPara = $d.Paragraphs.Last
$listPara.Style = "List Paragraph"
$listPara.Range.Text = "Make sure " + [char]0x201C + "related work" + [char]0x201D + " and other sections are tailored to the specific conference theme we are applying to."
$listPara.Range.ListFormat.ListId = 1
$listPara.Range.ListFormat.ListLevelNumber = 1
